# Update cryptocurrency market data snapshot (2024-02-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 51606
$ws.Range("E2").Value = 1013862656276
$ws.Range("F2").Value = 14670160193
$ws.Range("G2").Value = 1.0705
$ws.Range("D3").Value = 3032.71
$ws.Range("E3").Value = 364756261463
$ws.Range("F3").Value = 13335510805
$ws.Range("G3").Value = 2.53296
$ws.Range("D4").Value = 0.999706
$ws.Range("E4").Value = 97914782829
$ws.Range("F4").Value = 28324859926
$ws.Range("G4").Value = 0.00603
$ws.Range("D5").Value = 384.85
$ws.Range("E5").Value = 59145736339
$ws.Range("F5").Value = 884386510
$ws.Range("G5").Value = 1.21684
$ws.Range("D6").Value = 102.46
$ws.Range("E6").Value = 45211377309
$ws.Range("F6").Value = 1238311669
$ws.Range("G6").Value = 0.60151
$ws.Range("B7").Value = "STETH"
$ws.Range("C7").Value = "Lido Staked Ether"
$ws.Range("D7").Value = 3028.7
$ws.Range("E7").Value = 29707833002
$ws.Range("F7").Value = 22570638
$ws.Range("G7").Value = 2.45679
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "XRP"
$ws.Range("D8").Value = 0.543999
$ws.Range("E8").Value = 29706154673
$ws.Range("F8").Value = 589557514
$ws.Range("G8").Value = 0.00053
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "USDC"
$ws.Range("D9").Value = 0.999576
$ws.Range("E9").Value = 28171147681
$ws.Range("F9").Value = 3426957762
$ws.Range("G9").Value = -0.00023
$ws.Range("D10").Value = 0.588214
$ws.Range("E10").Value = 20690555656
$ws.Range("F10").Value = 327602582
$ws.Range("G10").Value = 0.45942
$ws.Range("D11").Value = 36.74
$ws.Range("E11").Value = 13862644871
$ws.Range("F11").Value = 360139432
$ws.Range("G11").Value = 0.34379
$ws.Range("D12").Value = 0.08595800000000001
$ws.Range("E12").Value = 12311390346
$ws.Range("F12").Value = 338438917
$ws.Range("G12").Value = 0.56513
$ws.Range("D13").Value = 0.13742
$ws.Range("E13").Value = 12099108728
$ws.Range("F13").Value = 228916962
$ws.Range("G13").Value = 0.00655
$ws.Range("D14").Value = 18.6
$ws.Range("E14").Value = 10907148471
$ws.Range("F14").Value = 273892365
$ws.Range("G14").Value = 1.57551
$ws.Range("D15").Value = 7.75
$ws.Range("E15").Value = 10348260060
$ws.Range("F15").Value = 183620894
$ws.Range("G15").Value = 0.05972
$ws.Range("D16").Value = 0.973957
$ws.Range("E16").Value = 9046271197
$ws.Range("F16").Value = 382779312
$ws.Range("G16").Value = -2.91932
$ws.Range("B17").Value = "UNI"
$ws.Range("C17").Value = "Uniswap"
$ws.Range("D17").Value = 10.77
$ws.Range("E17").Value = 8122875419
$ws.Range("F17").Value = 983919242
$ws.Range("G17").Value = -11.64245
$ws.Range("B18").Value = "WBTC"
$ws.Range("C18").Value = "Wrapped Bitcoin"
$ws.Range("D18").Value = 51610
$ws.Range("E18").Value = 8099961537
$ws.Range("F18").Value = 216902327
$ws.Range("G18").Value = 1.14176
$ws.Range("B19").Value = "TON"
$ws.Range("C19").Value = "Toncoin"
$ws.Range("D19").Value = 2.11
$ws.Range("E19").Value = 7320952881
$ws.Range("F19").Value = 18403467
$ws.Range("G19").Value = -0.96721
$ws.Range("B20").Value = "ICP"
$ws.Range("C20").Value = "Internet Computer"
$ws.Range("D20").Value = 12.43
$ws.Range("E20").Value = 5721594427
$ws.Range("F20").Value = 92311249
$ws.Range("G20").Value = 0.56468
$ws.Range("D21").Value = 0.00000961
$ws.Range("E21").Value = 5666839136
$ws.Range("F21").Value = 113235736
$ws.Range("G21").Value = 0.10224
$ws.Range("D22").Value = 266.96
$ws.Range("E22").Value = 5248773175
$ws.Range("F22").Value = 110840520
$ws.Range("G22").Value = -0.37352
$ws.Range("D23").Value = 70
$ws.Range("E23").Value = 5198847894
$ws.Range("F23").Value = 226573834
$ws.Range("G23").Value = 0.3166
$ws.Range("D24").Value = 0.999018
$ws.Range("E24").Value = 4950163204
$ws.Range("F24").Value = 111253911
$ws.Range("G24").Value = 0.07504
$ws.Range("B25").Value = "FIL"
$ws.Range("C25").Value = "Filecoin"
$ws.Range("D25").Value = 8.33
$ws.Range("E25").Value = 4288440075
$ws.Range("F25").Value = 449313043
$ws.Range("G25").Value = 5.22185
$ws.Range("B26").Value = "IMX"
$ws.Range("C26").Value = "Immutable"
$ws.Range("D26").Value = 3.08
$ws.Range("E26").Value = 4273190716
$ws.Range("F26").Value = 77724632
$ws.Range("G26").Value = -0.33811
$ws.Range("D27").Value = 10.28
$ws.Range("E27").Value = 3989900673
$ws.Range("F27").Value = 175946664
$ws.Range("G27").Value = -1.24392
$ws.Range("B28").Value = "KAS"
$ws.Range("C28").Value = "Kaspa"
$ws.Range("D28").Value = 0.172915
$ws.Range("E28").Value = 3931144496
$ws.Range("F28").Value = 33538391
$ws.Range("G28").Value = 3.85354
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "LEO Token"
$ws.Range("D29").Value = 4.21
$ws.Range("E29").Value = 3912524924
$ws.Range("F29").Value = 1309497
$ws.Range("G29").Value = 2.07278
$ws.Range("B30").Value = "NEAR"
$ws.Range("C30").Value = "NEAR Protocol"
$ws.Range("D30").Value = 3.69
$ws.Range("E30").Value = 3840784957
$ws.Range("F30").Value = 414525034
$ws.Range("G30").Value = 4.54301
$ws.Range("B31").Value = "ETC"
$ws.Range("C31").Value = "Ethereum Classic"
$ws.Range("D31").Value = 26.25
$ws.Range("E31").Value = 3762486384
$ws.Range("F31").Value = 110216642
$ws.Range("G31").Value = 1.68368
$ws.Range("B32").Value = "TAO"
$ws.Range("C32").Value = "Bittensor"
$ws.Range("D32").Value = 579.75
$ws.Range("E32").Value = 3659611425
$ws.Range("F32").Value = 14266389
$ws.Range("G32").Value = -2.55087
$ws.Range("B33").Value = "STX"
$ws.Range("C33").Value = "Stacks"
$ws.Range("D33").Value = 2.53
$ws.Range("E33").Value = 3654601547
$ws.Range("F33").Value = 70388012
$ws.Range("G33").Value = 2.27114
$ws.Range("B34").Value = "HBAR"
$ws.Range("C34").Value = "Hedera"
$ws.Range("D34").Value = 0.107599
$ws.Range("E34").Value = 3625964033
$ws.Range("F34").Value = 82402468
$ws.Range("G34").Value = -1.07704
$ws.Range("B35").Value = "OP"
$ws.Range("C35").Value = "Optimism"
$ws.Range("D35").Value = 3.7
$ws.Range("E35").Value = 3545709931
$ws.Range("F35").Value = 164549621
$ws.Range("G35").Value = 2.96799
$ws.Range("B36").Value = "APT"
$ws.Range("C36").Value = "Aptos"
$ws.Range("D36").Value = 9.539999999999999
$ws.Range("E36").Value = 3492554059
$ws.Range("F36").Value = 81133865
$ws.Range("G36").Value = 1.77096
$ws.Range("B37").Value = "XLM"
$ws.Range("C37").Value = "Stellar"
$ws.Range("D37").Value = 0.116357
$ws.Range("E37").Value = 3312669107
$ws.Range("F37").Value = 57358930
$ws.Range("G37").Value = -0.08545999999999999
$ws.Range("B38").Value = "FDUSD"
$ws.Range("C38").Value = "First Digital USD"
$ws.Range("D38").Value = 0.99766
$ws.Range("E38").Value = 3301453412
$ws.Range("F38").Value = 3016475587
$ws.Range("G38").Value = -0.35747
$ws.Range("B39").Value = "VET"
$ws.Range("C39").Value = "VeChain"
$ws.Range("D39").Value = 0.04457198
$ws.Range("E39").Value = 3242400882
$ws.Range("F39").Value = 78603193
$ws.Range("G39").Value = 2.56178
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "OKB"
$ws.Range("D40").Value = 50.5
$ws.Range("E40").Value = 3030010080
$ws.Range("F40").Value = 5659343
$ws.Range("G40").Value = -0.94948
$ws.Range("B41").Value = "INJ"
$ws.Range("C41").Value = "Injective"
$ws.Range("D41").Value = 33.95
$ws.Range("E41").Value = 3006829319
$ws.Range("F41").Value = 84559425
$ws.Range("G41").Value = -0.94517
$ws.Range("B42").Value = "LDO"
$ws.Range("C42").Value = "Lido DAO"
$ws.Range("D42").Value = 3.33
$ws.Range("E42").Value = 2967465062
$ws.Range("F42").Value = 89985202
$ws.Range("G42").Value = 1.47877
$ws.Range("B43").Value = "TIA"
$ws.Range("C43").Value = "Celestia"
$ws.Range("D43").Value = 17.01
$ws.Range("E43").Value = 2836340436
$ws.Range("F43").Value = 75099982
$ws.Range("G43").Value = 2.59609
$ws.Range("B44").Value = "RNDR"
$ws.Range("C44").Value = "Render"
$ws.Range("D44").Value = 7.47
$ws.Range("E44").Value = 2823864116
$ws.Range("F44").Value = 262530083
$ws.Range("G44").Value = 4.86404
$ws.Range("B45").Value = "GRT"
$ws.Range("C45").Value = "The Graph"
$ws.Range("D45").Value = 0.292679
$ws.Range("E45").Value = 2737733373
$ws.Range("F45").Value = 344194202
$ws.Range("G45").Value = 9.846679999999999
$ws.Range("B46").Value = "MNT"
$ws.Range("C46").Value = "Mantle"
$ws.Range("D46").Value = 0.797198
$ws.Range("E46").Value = 2573420780
$ws.Range("F46").Value = 64695810
$ws.Range("G46").Value = 2.81095
$ws.Range("D47").Value = 0.093056
$ws.Range("E47").Value = 2474964102
$ws.Range("F47").Value = 9614166
$ws.Range("G47").Value = -1.20863
$ws.Range("B48").Value = "ARB"
$ws.Range("C48").Value = "Arbitrum"
$ws.Range("D48").Value = 1.87
$ws.Range("E48").Value = 2380126193
$ws.Range("F48").Value = 239833081
$ws.Range("G48").Value = 1.65933
$ws.Range("B49").Value = "XMR"
$ws.Range("C49").Value = "Monero"
$ws.Range("D49").Value = 123.94
$ws.Range("E49").Value = 2250704723
$ws.Range("F49").Value = 33317637
$ws.Range("G49").Value = -0.91814
$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "Sei"
$ws.Range("D50").Value = 0.829756
$ws.Range("E50").Value = 2117336370
$ws.Range("F50").Value = 169098247
$ws.Range("G50").Value = 1.99686
$ws.Range("B51").Value = "SUI"
$ws.Range("C51").Value = "Sui"
$ws.Range("D51").Value = 1.64
$ws.Range("E51").Value = 1906957337
$ws.Range("F51").Value = 155310833
$ws.Range("G51").Value = 1.50402
